$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("dSF") values repulled/recalculated for several rows.
$updates = @{
    2  = -6
    3  = -1
    11 = 11
    12 = 1
    15 = -2
    16 = -3
    17 = -5
    19 = -5
    23 = 15
    24 = 4
    25 = 1
    28 = -4
    30 = 1
    31 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
